# Commit: "add functionality for table reading"
#
# C3 previously held =TRUE() (displayed with a custom "TRUE/FALSE" number
# format, numFmtId 165). Replace it with the literal formula 0, which also
# makes the old boolean display format (165) obsolete -> the cell now uses
# the plain "General" number format (164). Finally move the active
# selection from D8 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C3")

# Drop the boolean "TRUE"/"FALSE" custom format in favour of General,
# since the cell will no longer hold a boolean result.
$cell.NumberFormat = "General"

# Replace the =TRUE() formula with a plain 0.
$cell.Formula = "=0"

# Update the active selection to C3 (was D8).
$cell.Select()
